$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$style = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "22.039.91"
$ws.Range("D2").Style = $style
$ws.Range("E2").Value = "  -0.39%  "

$style = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.556.36"
$ws.Range("D3").Style = $style
$ws.Range("E3").Value = "  +0.25%  "

$style = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").Style = $style
$ws.Range("E4").Value = "  -0.18%  "

$style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9996"
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = "  -0.13%  "

$style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "290.73"
$ws.Range("D6").Style = $style
$ws.Range("E6").Value = "  +1.17%  "

$style = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3972"
$ws.Range("D7").Style = $style
$ws.Range("E7").Value = "  +3.93%  "

$style = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3226"
$ws.Range("D8").Style = $style
$ws.Range("E8").Value = "  -2.39%  "

$style = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "44.25"
$ws.Range("D9").Style = $style
$ws.Range("E9").Value = "  +1.09%  "

$style = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07252"
$ws.Range("D10").Style = $style
$ws.Range("E10").Value = "  -1.55%  "

$style = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.079"
$ws.Range("D11").Style = $style
$ws.Range("E11").Value = "  -4.80%  "

$style = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.000"
$ws.Range("D12").Style = $style
$ws.Range("E12").Value = "  -0.20%  "

$style = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.701"
$ws.Range("D13").Style = $style
$ws.Range("E13").Value = "  -1.96%  "

$style = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.81"
$ws.Range("D14").Style = $style
$ws.Range("E14").Value = "  -6.59%  "

$style = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001129"
$ws.Range("D15").Style = $style
$ws.Range("E15").Value = "  +5.34%  "

$style = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.640"
$ws.Range("D16").Style = $style
$ws.Range("E16").Value = "  -1.36%  "

$style = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.551.11"
$ws.Range("D17").Style = $style
$ws.Range("E17").Value = "  -0.21%  "

$style = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06600"
$ws.Range("D18").Style = $style
$ws.Range("E18").Value = "  -0.60%  "

$style = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "83.70"
$ws.Range("D19").Style = $style
$ws.Range("E19").Value = "  -2.60%  "

$ws.Range("E20").Value = "  -0.20%  "

$style = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.270"
$ws.Range("D21").Style = $style
$ws.Range("E21").Value = "  -1.40%  "

$style = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "15.57"
$ws.Range("D22").Style = $style
$ws.Range("E22").Value = "  -2.98%  "

$ws.Range("E23").Value = "  -3.14%  "

$style = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "22.059.48"
$ws.Range("D24").Style = $style
$ws.Range("E24").Value = "  -0.28%  "

$style = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.364"
$ws.Range("D25").Style = $style
$ws.Range("E25").Value = "  +3.36%  "

$style = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.422"
$ws.Range("D26").Style = $style
$ws.Range("E26").Value = "  -4.49%  "

$style = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "148.67"
$ws.Range("D27").Style = $style
$ws.Range("E27").Value = "  -1.20%  "

$style = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.62"
$ws.Range("D28").Style = $style
$ws.Range("E28").Value = "  -2.68%  "

$style = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.871"
$ws.Range("D29").Style = $style
$ws.Range("E29").Value = "  -0.95%  "

$style = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.727.32"
$ws.Range("D30").Style = $style
$ws.Range("E30").Value = "  -0.78%  "

$style = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "119.32"
$ws.Range("D31").Style = $style
$ws.Range("E31").Value = "  -2.38%  "

$style = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9937"
$ws.Range("D32").Style = $style
$ws.Range("E32").Value = "  -8.60%  "

$style = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.843"
$ws.Range("D33").Style = $style
$ws.Range("E33").Value = "  -0.54%  "

$style = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08319"
$ws.Range("D34").Style = $style
$ws.Range("E34").Value = "  +1.43%  "

$style = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.179"
$ws.Range("D35").Style = $style
$ws.Range("E35").Value = "  -1.25%  "

$style = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.602"
$ws.Range("D36").Style = $style
$ws.Range("E36").Value = "  -15.55%  "

$style = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02264"
$ws.Range("D37").Style = $style
$ws.Range("E37").Value = "  -2.38%  "

$style = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.130"
$ws.Range("D38").Style = $style
$ws.Range("E38").Value = "  -3.23%  "

$style = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06006"
$ws.Range("D39").Style = $style
$ws.Range("E39").Value = "  -4.48%  "

$style = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.209"
$ws.Range("D40").Style = $style
$ws.Range("E40").Value = "  -2.22%  "

$style = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2042"
$ws.Range("D41").Style = $style
$ws.Range("E41").Value = "  -4.92%  "

$style = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.001"
$ws.Range("D42").Style = $style
$ws.Range("E42").Value = "  -0.03%  "

$style = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.77"
$ws.Range("D43").Style = $style
$ws.Range("E43").Value = "  -1.91%  "

$style = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5832"
$ws.Range("D44").Style = $style
$ws.Range("E44").Value = "  -3.55%  "

$style = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.09"
$ws.Range("D45").Style = $style
$ws.Range("E45").Value = "  -5.15%  "

$ws.Range("E46").Value = "  +0.45%  "

$style = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5581"
$ws.Range("D47").Style = $style
$ws.Range("E47").Value = "  -4.49%  "

$style = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "118.39"
$ws.Range("D48").Style = $style
$ws.Range("E48").Value = "  -2.74%  "

$style = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.902"
$ws.Range("D49").Style = $style
$ws.Range("E49").Value = "  -3.33%  "

$ws.Range("E50").Value = "  -2.90%  "

$style = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06830"
$ws.Range("D51").Style = $style
$ws.Range("E51").Value = "  -2.78%  "
